$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new "Battery charger" line under the Battery rows (new row 13) ---
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "Battery charger"
$ws.Range("B13").Value = "IMAX B6AC V2 Professional Balance Charger/Discharger"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 38.19
$ws.Range("E13").Formula = "=C13*D13"

# --- 2. Transmission section: bump the 12-tooth sprocket quantity to 2 (now row 16) ---
$ws.Range("C16").Value = 2

# --- 3. Roller chain row becomes "...DIN ISO 606 - 1m" with quantity 5 (now row 18) ---
$ws.Range("B18").Value = "3/8“ x 7/32“ DIN ISO 606 - 1m"
$ws.Range("C18").Value = 5

# --- 4. Insert a new "Connecting link" row right after the roller chain row (new row 19) ---
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "Connecting link"
$ws.Range("B19").Value = "3/8“ x 7/32“ DIN 8187"
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 0
$ws.Range("E19").Formula = "=C19*D19"

# --- 5. Taperlock quantity drops from 2 to 1 (now row 20) ---
$ws.Range("C20").Value = 1

# --- 6. Case/box: update the waterproof ABS box dimensions (now row 31) ---
$ws.Range("B31").Value = "Waterproof ABS box - Dark grey - 195 x 80 x 55mm G353"
